$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 33, contains a date value of 46074 that must be
# advanced by one day to 46075 (2026-02-21 -> 2026-02-22) for every row.
$ws.Range("C2:C33").Value = 46075
